# [ADD] New normalize way
# Update GlobalPriority (column B) and NivelSeguridad (column C) values
# for each Alternative row with the newly computed normalized figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.109029252404279
$ws.Range("C2").Value = 2

$ws.Range("B3").Value = 0.07476217823650765

$ws.Range("B4").Value = 0.176267729111987

$ws.Range("B5").Value = 0.09612066273305153
$ws.Range("C5").Value = 2

$ws.Range("B6").Value = 0.1691909686138251

$ws.Range("B7").Value = 0.1777097446305805

$ws.Range("B8").Value = 0.09585285992210585
$ws.Range("C8").Value = 2

$ws.Range("B9").Value = 0.1473826792247265
$ws.Range("C9").Value = 2
